$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.06750333582058965
$ws.Range("C2").Value = 0.9987668739163376
$ws.Range("B3").Value = 0.1259519652509079
$ws.Range("C3").Value = 0.9906942099676707
